# Update the workbook to add data for 2022-09-06
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the header label to reflect the new "through" date
$ws.Name = "Through 2022-08-29"
$ws.Range("I1").Value = "2022 (through 08-29)"

# Update the two changed data values
$ws.Range("I9").Value = 152
$ws.Range("I14").Value = 1123
